# "tests du générateur d'onde triangulaire, évantuelle modification du module de test"
#
# The workbook's last sheet ("test_adc_dac_modele_rc") holds a small table of
# test-tension measurements (columns B..E, rows 3..10). This edit:
#   - drops the unused "tension test2" header in column E (column E becomes
#     completely empty / out of the sheet's used range),
#   - fills in the (until now empty) "?" placeholder measurements in column C
#     for every data row (4..10), pending the real values from the triangular
#     wave generator test,
#   - leaves the cursor/selection on C14, where the user ended up after typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# Remove the now-unused "tension test2" column header.
$ws.Range("E3").ClearContents()

# Fill in the new placeholder measurement column.
$ws.Range("C4").Value = "?"
$ws.Range("C5").Value = "?"
$ws.Range("C6").Value = "?"
$ws.Range("C7").Value = "?"
$ws.Range("C8").Value = "?"
$ws.Range("C9").Value = "?"
$ws.Range("C10").Value = "?"

# Match the author's final cursor position on that sheet.
$ws.Range("C14").Select() | Out-Null
